$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reviewer Table 1")

$ws.Range("B3").Value = "BLCA: 35 (90)`nBLSC: 2 (5.1)`nUCU: 2 (5.1)`ncomplete: n = 39"
$ws.Range("B4").Value = "66 (57, 80)`ncomplete: n = 24"
$ws.Range("B5").Value = "Female: 15 (38)`nMale: 22 (56)`nUnknown: 2 (5.1)`ncomplete: n = 39"
$ws.Range("B6").Value = "7 (18)`ncomplete: n = 39"

$ws.Range("B7").ClearContents()
$ws.Range("B8").ClearContents()
$ws.Range("B9").ClearContents()
